$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row of data: year 2019, index value 0.1
$ws.Range("A26").Value = 2019
$ws.Range("B26").Value = 0.1

# Update selection to reflect the new active cell after data entry
$ws.Range("A27").Select()
